# Generate Report for Handoff
# Adds a new row (row 3) for file e9a51b8e-d50e-45ec-85d9-7b63b723ef86
# to all three sheets: Overview, zh-cn, de-de.

$wb = $excel.ActiveWorkbook

$commit = "9f4a8a16461f5a29c907229f23b974446059a1d2"
$hyperColor = 15570276   # VBA RGB() value equivalent to rgb="FF6495ED"

# ---------------------------------------------------------------
# Sheet 1: Overview
# ---------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$tOverview = $wsOverview.ListObjects.Item(1)
$tOverview.ListRows.Add() | Out-Null

$wsOverview.Range("A3").Value = "e9a51b8e-d50e-45ec-85d9-7b63b723ef86.md"
$wsOverview.Range("B3").Value = "e2e\e9a51b8e-d50e-45ec-85d9-7b63b723ef86.md"
$wsOverview.Range("C3").Value = ".md"
$wsOverview.Range("D3").Value = "'"
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-08-18 16:38:47"
$wsOverview.Range("G3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/$commit/e2e/e9a51b8e-d50e-45ec-85d9-7b63b723ef86.md", $null, $null, "e2e\e9a51b8e-d50e-45ec-85d9-7b63b723ef86.md") | Out-Null
$wsOverview.Range("B3").Font.Underline = $true
$wsOverview.Range("B3").Font.Color = $hyperColor

# ---------------------------------------------------------------
# Sheet 2: zh-cn
# ---------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")
$tZh = $wsZh.ListObjects.Item(1)
$tZh.ListRows.Add() | Out-Null

$wsZh.Range("A3").Value = "e9a51b8e-d50e-45ec-85d9-7b63b723ef86.md"
$wsZh.Range("B3").Value = ".md"
$wsZh.Range("C3").Value = "Ready for handoff"
$wsZh.Range("D3").Value = "e2e"
$wsZh.Range("E3").Value = "ht"
$wsZh.Range("F3").Value = "'False"
$wsZh.Range("G3").Value = "e9a51b8e-d50e-45ec-85d9-7b63b723ef86.5d803a5fba518d35c82b52376a7a2194ce50e463.zh-cn.xlf"
$wsZh.Range("H3").Value = "2016-08-18 16:38:42"
$wsZh.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("I3").Value = "'"
$wsZh.Range("J3").Value = "'"
$wsZh.Range("K3").Value = "0001-01-01 00:00:00"
$wsZh.Range("K3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("L3").Value = "'"
$wsZh.Range("M3").Value = "'True"
$wsZh.Range("N3").Value = "'"
$wsZh.Range("O3").Value = "'False"
$wsZh.Range("P3").Value = "'"

$wsZh.Hyperlinks.Add($wsZh.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/$commit/e2e/e9a51b8e-d50e-45ec-85d9-7b63b723ef86.md", $null, $null, "e9a51b8e-d50e-45ec-85d9-7b63b723ef86.md") | Out-Null
$wsZh.Range("A3").Font.Underline = $true
$wsZh.Range("A3").Font.Color = $hyperColor

# ---------------------------------------------------------------
# Sheet 3: de-de
# ---------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")
$tDe = $wsDe.ListObjects.Item(1)
$tDe.ListRows.Add() | Out-Null

$wsDe.Range("A3").Value = "e9a51b8e-d50e-45ec-85d9-7b63b723ef86.md"
$wsDe.Range("B3").Value = ".md"
$wsDe.Range("C3").Value = "Ready for handoff"
$wsDe.Range("D3").Value = "e2e"
$wsDe.Range("E3").Value = "ht"
$wsDe.Range("F3").Value = "'False"
$wsDe.Range("G3").Value = "e9a51b8e-d50e-45ec-85d9-7b63b723ef86.5d803a5fba518d35c82b52376a7a2194ce50e463.de-de.xlf"
$wsDe.Range("H3").Value = "2016-08-18 16:38:47"
$wsDe.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("K3").Value = "0001-01-01 00:00:00"
$wsDe.Range("K3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("M3").Value = "'True"
$wsDe.Range("O3").Value = "'False"

$wsDe.Hyperlinks.Add($wsDe.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/$commit/e2e/e9a51b8e-d50e-45ec-85d9-7b63b723ef86.md", $null, $null, "e9a51b8e-d50e-45ec-85d9-7b63b723ef86.md") | Out-Null
$wsDe.Range("A3").Font.Underline = $true
$wsDe.Range("A3").Font.Color = $hyperColor
